# Update "想去人数" (F column) figures across the workbook's sheets,
# matching the refreshed bilibili show-interest export.
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibitions)
$ws = $wb.Worksheets.Item(1)
$ws.Range("F4").Value = 415
$ws.Range("F5").Value = 1242
$ws.Range("F7").Value = 7501
$ws.Range("F9").Value = 106
$ws.Range("F10").Value = 2075
$ws.Range("F11").Value = 8169
$ws.Range("F14").Value = 5576
$ws.Range("F16").Value = 2536
$ws.Range("F17").Value = 1088
$ws.Range("F18").Value = 4573
$ws.Range("F22").Value = 24
$ws.Range("F23").Value = 454
$ws.Range("F24").Value = 1749
$ws.Range("F25").Value = 30
$ws.Range("F26").Value = 2698
$ws.Range("F28").Value = 308
$ws.Range("F29").Value = 107
$ws.Range("F30").Value = 247
$ws.Range("F31").Value = 620
$ws.Range("F33").Value = 528
$ws.Range("F34").Value = 1602
$ws.Range("F36").Value = 13
$ws.Range("F37").Value = 2549
$ws.Range("F38").Value = 2261
$ws.Range("F40").Value = 22
$ws.Range("F41").Value = 318

# Sheet 2: 演出 (Performances)
$ws = $wb.Worksheets.Item(2)
$ws.Range("F6").Value = 31

# Sheet 3: 本地生活 (Local life)
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 260
$ws.Range("F3").Value = 1291

# Sheet 4: 全部类型 (All types, combined listing)
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 260
$ws.Range("F4").Value = 1291
$ws.Range("F6").Value = 415
$ws.Range("F7").Value = 1242
$ws.Range("F9").Value = 7501
$ws.Range("F11").Value = 106
$ws.Range("F12").Value = 2075
$ws.Range("F13").Value = 8169
$ws.Range("F16").Value = 5576
$ws.Range("F18").Value = 2536
$ws.Range("F19").Value = 1088
$ws.Range("F20").Value = 4573
$ws.Range("F25").Value = 24
$ws.Range("F27").Value = 454
$ws.Range("F28").Value = 1750
$ws.Range("F29").Value = 30
$ws.Range("F30").Value = 2698
$ws.Range("F32").Value = 308
$ws.Range("F33").Value = 107
$ws.Range("F34").Value = 247
$ws.Range("F36").Value = 620
$ws.Range("F38").Value = 528
$ws.Range("F40").Value = 1602
$ws.Range("F42").Value = 13
$ws.Range("F43").Value = 2549
$ws.Range("F44").Value = 31
$ws.Range("F45").Value = 2261
$ws.Range("F47").Value = 22
$ws.Range("F48").Value = 318

Write-Output "Applied all F-column updates"
